$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the first two data rows (original rows 2 and 3).
# This shifts the old rows 4 and 5 up to become the new rows 2 and 3,
# matching the target diff.
$ws.Rows.Item(2).Delete()
$ws.Rows.Item(2).Delete()
